$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "97.043.97"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
$ws.Range("D3").Value = "3.694.24"
$ws.Range("E3").Value = "  +2.97%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'239.97"
$ws.Range("E5").Value = "  -0.80%  "

# Row 6
$ws.Range("D6").Value = "'1.89"
$ws.Range("E6").Value = "  +8.84%  "

# Row 7
$ws.Range("D7").Value = "'654.67"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  -1.41%  "

# Row 9
$ws.Range("E9").Value = "  +2.71%  "

# Row 10
$ws.Range("D10").Value = "'0.999"
$ws.Range("E10").Value = "  +0.01%  "

# Row 11
$ws.Range("D11").Value = "3.694.57"
$ws.Range("E11").Value = "  +3.07%  "

# Row 12
$ws.Range("D12").Value = "'45.54"
$ws.Range("E12").Value = "  +2.38%  "

# Row 13
$ws.Range("E13").Value = "  +0.91%  "

# Row 14
$ws.Range("D14").Value = "'6.90"
$ws.Range("E14").Value = "  +6.69%  "

# Row 15
$ws.Range("D15").Value = "4.380.50"
$ws.Range("E15").Value = "  +2.98%  "

# Row 16
$ws.Range("D16").Value = "'0.0000268"
$ws.Range("E16").Value = "  +2.50%  "

# Row 17
$ws.Range("D17").Value = "96.905.58"
$ws.Range("E17").Value = "  -0.22%  "

# Row 18
$ws.Range("D18").Value = "'9.06"
$ws.Range("E18").Value = "  +3.38%  "

# Row 19
$ws.Range("D19").Value = "3.701.35"
$ws.Range("E19").Value = "  +3.27%  "

# Row 20
$ws.Range("D20").Value = "'19.40"
$ws.Range("E20").Value = "  +6.47%  "

# Row 21
$ws.Range("D21").Value = "'12.88"
$ws.Range("E21").Value = "  +1.84%  "

# Row 22
$ws.Range("D22").Value = "'0.529"
$ws.Range("E22").Value = "  -0.63%  "

# Row 23
$ws.Range("D23").Value = "'527.36"
$ws.Range("E23").Value = "  +1.66%  "

# Row 24
$ws.Range("D24").Value = "'3.52"
$ws.Range("E24").Value = "  +0.49%  "

# Row 25
$ws.Range("D25").Value = "'7.16"
$ws.Range("E25").Value = "  +2.73%  "

# Row 26
$ws.Range("D26").Value = "'0.0000205"
$ws.Range("E26").Value = "  -1.39%  "

# Row 27
$ws.Range("D27").Value = "'102.41"
$ws.Range("E27").Value = "  +0.35%  "

# Row 28
$ws.Range("D28").Value = "'13.42"
$ws.Range("E28").Value = "  +1.94%  "

# Row 29
$ws.Range("D29").Value = "'0.169"
$ws.Range("E29").Value = "  -1.73%  "

# Row 30
$ws.Range("D30").Value = "'12.55"
$ws.Range("E30").Value = "  +3.95%  "

# Row 31
$ws.Range("D31").Value = "'3.05"
$ws.Range("E31").Value = "  +1.78%  "

# Row 32
$ws.Range("E32").Value = "  -0.06%  "

# Row 33
$ws.Range("D33").Value = "'1.88"
$ws.Range("E33").Value = "  +14.76%  "

# Row 34
$ws.Range("D34").Value = "'0.187"
$ws.Range("E34").Value = "  +0.64%  "

# Row 35
$ws.Range("D35").Value = "'32.81"
$ws.Range("E35").Value = "  +2.35%  "

# Row 36
$ws.Range("E36").Value = "  +0.12%  "

# Row 37
$ws.Range("D37").Value = "'0.613"
$ws.Range("E37").Value = "  +6.81%  "

# Row 38
$ws.Range("D38").Value = "'658.12"
$ws.Range("E38").Value = "  +6.65%  "

# Row 39
$ws.Range("D39").Value = "'9.09"
$ws.Range("E39").Value = "  +3.39%  "

# Row 40
$ws.Range("D40").Value = "'7.00"
$ws.Range("E40").Value = "  +15.70%  "

# Row 41
$ws.Range("D41").Value = "'0.164"
$ws.Range("E41").Value = "  +5.87%  "

# Row 42
$ws.Range("D42").Value = "'2.02"
$ws.Range("E42").Value = "  +3.54%  "

# Row 43
$ws.Range("D43").Value = "'0.969"
$ws.Range("E43").Value = "  +4.00%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'38.69"
$ws.Range("E44").Value = "  +17.14%  "

# Row 45
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.03%  "

# Row 46
$ws.Range("D46").Value = "'0.453"
$ws.Range("E46").Value = "  +7.70%  "

# Row 47
$ws.Range("D47").Value = "'0.0462"

# Row 48
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.78"
$ws.Range("E49").Value = "  +2.16%  "

# Row 50
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'23.64"
$ws.Range("E50").Value = "  -0.05%  "

# Row 51
$ws.Range("D51").Value = "'3.60"
$ws.Range("E51").Value = "  +2.91%  "
